$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.179.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.531.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.531.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.993.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.170.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.555.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.671.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.356"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.562"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
